$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("participants")

# Rename the "gestational age (weeks)" header to "gestational_age_(weeks)"
$ws.Range("E1").Value2 = "gestational_age_(weeks)"

# Add the new "pregnancy_term" column header
$ws.Range("L1").Value2 = "pregnancy_term"

# Update the race list-validation to add "More than one" / "Unknown or not reported"
# in place of "Multiple", keeping its original position among the validations.
$raceValidation = $ws.Range("G2:G1001").Validation
$raceValidation.Modify(3, 1, 1, '"American Indian or Alaska Native,Asian,Native Hawaiian or Other Pacific Islander,Black or African American,White,More than one,Unknown or not reported"')
